$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the current last row (row 83), shifting nothing below it
# (row 84 doesn't exist yet) and copy the existing row 83 values/format into it.
$ws.Rows.Item(84).Insert()

# Copy row 83 (the old data, about to be superseded) down into the new row 84.
$ws.Rows.Item(83).Copy()
$ws.Rows.Item(84).PasteSpecial()

# Now update row 83 in place with the newer weekly figures.
$ws.Cells.Item(83, 4).Value = 45239
$ws.Cells.Item(83, 11).Value = 2000
$ws.Cells.Item(83, 12).Value = 2000
$ws.Cells.Item(83, 13).Value = 2000
$ws.Cells.Item(83, 16).Value = 2000
